$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -5.889099999999998
$ws.Range("C9").Value = -11.8376
$ws.Range("C18").Value = -14.2959
$ws.Range("C20").Value = -13.69369999999998
$ws.Range("E21").Value = 13.0719
